$d = $word.ActiveDocument

# The document currently ends with a single paragraph "12.04.2017" that
# carries the (hidden) "_GoBack" bookmark right after its text. We need to
# append two new paragraphs after it, with the bookmark ending up after the
# text of the final new paragraph (mirroring the diff).
#
# Inserting text that already contains paragraph marks directly "before"
# the bookmark's range tends to leave the bookmark anchored in the wrong
# spot (before the run instead of after it). Instead we insert the new
# text as a single run using a placeholder separator, then use Find &
# Replace to turn the separator into real paragraph breaks — this keeps
# the bookmark correctly positioned after the final run.

$bm = $d.Bookmarks("_GoBack")
$sep = "~~NEWPARA~~"

$bm.Range.InsertBefore($sep + "- getestet wie man auf einer Website mehrere Ansichten hat (PHP include)" + $sep + "-nächster Schritt – schauen ob die includierungen RealTime sind")

$d.Content.Find.Execute($sep, $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)
